$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '256.74'
Set-TextCell $ws 'E2' '-1.39%'
Set-TextCell $ws 'G2' '18'
Set-TextCell $ws 'D3' '27.34'
Set-TextCell $ws 'E3' '-1.95%'
Set-TextCell $ws 'G3' '18'
Set-TextCell $ws 'D4' '4.536'
Set-TextCell $ws 'E4' '-13.16%'
Set-TextCell $ws 'G4' '18'
Set-TextCell $ws 'D5' '0.05887'
Set-TextCell $ws 'E5' '-0.65%'
Set-TextCell $ws 'G5' '18'
Set-TextCell $ws 'E6' '-1.53%'
Set-TextCell $ws 'G6' '18'
Set-TextCell $ws 'D7' '0.8574'
Set-TextCell $ws 'E7' '-1.86%'
Set-TextCell $ws 'G7' '18'
Set-TextCell $ws 'D8' '0.9258'
Set-TextCell $ws 'E8' '-7.87%'
Set-TextCell $ws 'G8' '18'
Set-TextCell $ws 'B9' 'One'
Set-TextCell $ws 'C9' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell $ws 'D9' '0.01038'
Set-TextCell $ws 'E9' '1,612.03%'
Set-TextCell $ws 'G9' '18'
Set-TextCell $ws 'B10' 'WazirX'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws 'D10' '0.1406'
Set-TextCell $ws 'E10' '-1.36%'
Set-TextCell $ws 'G10' '18'
Set-TextCell $ws 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws 'D11' '0.03635'
Set-TextCell $ws 'E11' '-0.02%'
Set-TextCell $ws 'G11' '18'
Set-TextCell $ws 'B12' 'MandalaExchangeToken'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell $ws 'D12' '0.07087'
Set-TextCell $ws 'E12' '-1.98%'
Set-TextCell $ws 'G12' '18'
Set-TextCell $ws 'B13' 'BitrueCoin'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell $ws 'D13' '0.03232'
Set-TextCell $ws 'E13' '1.09%'
Set-TextCell $ws 'G13' '18'
Set-TextCell $ws 'B14' 'BitMartToken'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell $ws 'D14' '0.09206'
Set-TextCell $ws 'E14' '-0.37%'
Set-TextCell $ws 'G14' '18'
Set-TextCell $ws 'B15' 'BitForexToken'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell $ws 'D15' '0.001542'
Set-TextCell $ws 'E15' '0.12%'
Set-TextCell $ws 'G15' '18'
Set-TextCell $ws 'D16' '0.006059'
Set-TextCell $ws 'E16' '4.01%'
Set-TextCell $ws 'G16' '18'
Set-TextCell $ws 'D17' '3.518'
Set-TextCell $ws 'E17' '0.76%'
Set-TextCell $ws 'G17' '18'
Set-TextCell $ws 'D18' '3.198'
Set-TextCell $ws 'E18' '-1.24%'
Set-TextCell $ws 'G18' '18'
Set-TextCell $ws 'E19' '-0.77%'
Set-TextCell $ws 'G19' '18'
Set-TextCell $ws 'D20' '0.3058'
Set-TextCell $ws 'E20' '-2.05%'
Set-TextCell $ws 'G20' '18'
Set-TextCell $ws 'E21' '-1.01%'
Set-TextCell $ws 'G21' '18'
Set-TextCell $ws 'D22' '3.847'
Set-TextCell $ws 'E22' '9.17%'
Set-TextCell $ws 'G22' '18'
Set-TextCell $ws 'D23' '0.04209'
Set-TextCell $ws 'G23' '18'
Set-TextCell $ws 'D24' '0.001218'
Set-TextCell $ws 'E24' '0.10%'
Set-TextCell $ws 'G24' '18'
Set-TextCell $ws 'D25' '0.004279'
Set-TextCell $ws 'E25' '-6.34%'
Set-TextCell $ws 'G25' '18'
Set-TextCell $ws 'E26' '0.13%'
Set-TextCell $ws 'G26' '18'
Set-TextCell $ws 'D27' '0.0001510'
Set-TextCell $ws 'E27' '-21.96%'
Set-TextCell $ws 'G27' '18'
Set-TextCell $ws 'G28' '18'
Set-TextCell $ws 'G29' '18'
Set-TextCell $ws 'G30' '18'
Set-TextCell $ws 'G31' '18'
Set-TextCell $ws 'G32' '18'
Set-TextCell $ws 'G33' '18'
Set-TextCell $ws 'G34' '18'
Set-TextCell $ws 'G35' '18'
Set-TextCell $ws 'G36' '18'
Set-TextCell $ws 'G37' '18'
Set-TextCell $ws 'G38' '18'
Set-TextCell $ws 'G39' '18'
Set-TextCell $ws 'D40' '0.03831'
Set-TextCell $ws 'E40' '-0.70%'
Set-TextCell $ws 'G40' '18'
Set-TextCell $ws 'B41' 'KickToken'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell $ws 'D41' '0.006240'
Set-TextCell $ws 'E41' '13.60%'
Set-TextCell $ws 'G41' '18'
Set-TextCell $ws 'B42' 'BKEXToken'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws 'D42' '0.1098'
Set-TextCell $ws 'E42' '-1.09%'
Set-TextCell $ws 'G42' '18'
Set-TextCell $ws 'D43' '0.002200'
Set-TextCell $ws 'E43' '-7.31%'
Set-TextCell $ws 'G43' '18'
Set-TextCell $ws 'E44' '3.97%'
Set-TextCell $ws 'G44' '18'
Set-TextCell $ws 'D45' '0.00005470'
Set-TextCell $ws 'E45' '0.86%'
Set-TextCell $ws 'G45' '18'
Set-TextCell $ws 'E46' '0.13%'
Set-TextCell $ws 'G46' '18'
Set-TextCell $ws 'D47' '0.3000'
Set-TextCell $ws 'E47' '251.32%'
Set-TextCell $ws 'G47' '18'
Set-TextCell $ws 'D48' '0.1054'
Set-TextCell $ws 'E48' '4,831.55%'
Set-TextCell $ws 'G48' '18'
Set-TextCell $ws 'E49' '0.13%'
Set-TextCell $ws 'G49' '18'
Set-TextCell $ws 'E50' '0.13%'
Set-TextCell $ws 'G50' '18'
Set-TextCell $ws 'G51' '18'

Write-Output "Applied 133 cell updates"
